$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 0.58486341587430513
$ws.Range("C2").Value = 1.3708319777246605
$ws.Range("D2").Value = 3.2922372724995066
$ws.Range("E2").Value = 0.99785890624759799

$ws.Range("B3").Value = 1.8997263969175724
$ws.Range("C3").Value = 0.49570421829894007
$ws.Range("D3").Value = 1.6867536660769435
$ws.Range("E3").Value = 0.19224433899932919

$ws.Range("B1:E3").Select()
